$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 90. This shifts the old row 90 down to row 91,
# and leaves an empty (but style-inherited) row 90 ready to receive the data
# that used to live in the lower half of old row 89.
$ws.Rows.Item(90).EntireRow.Insert()

# --- Row 88: update values (date, volume, prices) ---
$ws.Range("D88").Value = 44747
$ws.Range("J88").Value = 300
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = 3750
$ws.Range("P88").Value = 1875

# --- Row 89: update values (date, volume, prices, unit text) ---
$ws.Range("D89").Value = 44566
$ws.Range("J89").Value = 250
$ws.Range("K89").Value = 4000
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 4250
$ws.Range("N89").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("P89").Value = 2125
$ws.Range("Q89").Value = 2

# --- Row 90 (new row): fill in with the data that used to be in old row 89 ---
$ws.Range("A90").Value = 1
$ws.Range("B90").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C90").Value = "Arica y Parinacota"
$ws.Range("D90").Value = 44637
$ws.Range("E90").Value = 15
$ws.Range("F90").Value = 100112038
$ws.Range("G90").Value = "Cebollín baby"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 300
$ws.Range("K90").Value = 1400
$ws.Range("L90").Value = 1500
$ws.Range("M90").Value = 1450
$ws.Range("N90").Value = "$/paquete"
$ws.Range("O90").Value = "Región de Arica y Parinacota"
$ws.Range("P90").Value = 1450
$ws.Range("Q90").Value = 1
$ws.Range("R90").Value = "Hortaliza"
